# Added Password Encryption code
# - Renames test-case style labels to include TC01_/TC02_ prefixes
# - Encrypts (base64-encodes) the plaintext passwords stored in the Login sheet
# - Switches the active/selected sheet from Admin back to Login

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Login" sheet
$ws2 = $wb.Worksheets.Item(2)   # "Admin" sheet

function Protect-Password([string]$plainText) {
    $bytes = [System.Text.Encoding]::UTF8.GetBytes($plainText)
    return [Convert]::ToBase64String($bytes)
}

# --- Login sheet: rename test case identifiers ---
$ws1.Range("A1").Value = "TC01_" + $ws1.Range("A1").Value()
$ws1.Range("D2").Value = $ws1.Range("A1").Value()

$ws1.Range("A3").Value = "TC02_" + $ws1.Range("A3").Value()
$ws1.Range("E6").Value = $ws1.Range("A3").Value()

# --- Login sheet: encrypt (base64 encode) stored passwords ---
$ws1.Range("C2").Value = Protect-Password($ws1.Range("C2").Value())
$ws1.Range("C5").Value = Protect-Password($ws1.Range("C5").Value())

$ws1.Range("C4").Value = Protect-Password($ws1.Range("C4").Value())
$ws1.Range("C6").Value = Protect-Password($ws1.Range("C6").Value())

# --- Admin sheet: rename test case identifier ---
$ws2.Range("A1").Value = "TC01_Admin_" + $ws2.Range("A1").Value()
$ws2.Range("E2").Value = $ws2.Range("A1").Value()

# --- Switch active sheet/selection back to Login, away from Admin ---
$ws2.Range("E8").Select() | Out-Null
$ws1.Activate() | Out-Null
